$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 301.070175
$ws.Range("D2").Value = 127.794305

$ws.Range("B3").Value = 8.451268000000001
$ws.Range("D3").Value = 1.793641
$ws.Range("E3").Value = 0.167962

$ws.Range("B4").Value = 782.157687
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = 0.337557
$ws.Range("H5").Value = -0.218939
$ws.Range("I5").Value = 0.894053
$ws.Range("J5").Value = 0.327635

$ws.Range("G6").Value = 0.468996
$ws.Range("H6").Value = -0.115514
$ws.Range("I6").Value = 1.053505
$ws.Range("J6").Value = 0.143438

$ws.Range("G7").Value = 0.131439
$ws.Range("H7").Value = -0.308073
$ws.Range("I7").Value = 0.570951
$ws.Range("J7").Value = 0.761257
